{"js": "// Corrections to UNICEF sitrep (November / December) \"Access to education\" figures.\n// Each entry is [oldText, newText]; all values are unique in the document, so a\n// direct search-and-replace on each one is safe and unambiguous.\nconst replacements = [\n  [\"170,873\", \"152,186\"],\n  [\"131,473\", \"147,338\"],\n  [\"19.72\", \"17.56\"],\n  [\"450,876\", \"461,832\"],\n  [\"407,635\", \"456,984\"],\n  [\"37.54\", \"38.45\"],\n  [\"74,534\", \"68,300\"],\n  [\"58,790\", \"66,480\"],\n  [\"17.20\", \"15.76\"],\n  [\"198,668\", \"205,488\"],\n  [\"181,035\", \"203,668\"],\n  [\"33.08\", \"34.22\"],\n  [\"96,339\", \"83,886\"],\n  [\"72,683\", \"80,858\"],\n  [\"22.23\", \"19.36\"],\n  [\"252,208\", \"256,344\"],\n  [\"226,600\", \"253,316\"],\n  [\"41.99\", \"42.68\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Corrections to UNICEF sitrep (November / December) \"Access to education\" figures.\n# Each pair is (oldText, newText); all values are unique in the document, so a\n# direct Find/Replace on each one is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"170,873\", \"152,186\"),\n    @(\"131,473\", \"147,338\"),\n    @(\"19.72\",   \"17.56\"),\n    @(\"450,876\", \"461,832\"),\n    @(\"407,635\", \"456,984\"),\n    @(\"37.54\",   \"38.45\"),\n    @(\"74,534\",  \"68,300\"),\n    @(\"58,790\",  \"66,480\"),\n    @(\"17.20\",   \"15.76\"),\n    @(\"198,668\", \"205,488\"),\n    @(\"181,035\", \"203,668\"),\n    @(\"33.08\",   \"34.22\"),\n    @(\"96,339\",  \"83,886\"),\n    @(\"72,683\",  \"80,858\"),\n    @(\"22.23\",   \"19.36\"),\n    @(\"252,208\", \"256,344\"),\n    @(\"226,600\", \"253,316\"),\n    @(\"41.99\",   \"42.68\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"No match found for '$oldText'\"\n    }\n}\n"}
